# Update documentation to include SQLs
#
# After the paragraph "This is because if we put createDatabase and
# useDatabase inside a same SQL, when we execute it at runtime JDBC will
# perform both queries and throw and error." insert four new list
# paragraphs describing the SQLs / insertTaiKhoan query.

$d = $word.ActiveDocument

# Locate the anchor paragraph via Find.
$findRange = $d.Content
$findRange.Find.Execute( `
    "This is because if we put createDatabase and useDatabase inside a same SQL, when we execute it at runtime JDBC will perform both queries and throw and error.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$anchorParagraph = $findRange.Paragraphs(1)
$anchorRange = $anchorParagraph.Range

# Determine the 1-based paragraph index of the anchor paragraph so we can
# walk forward by index as new paragraphs get inserted after it.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $anchorRange.Start) {
        $anchorIndex = $i
        break
    }
}

# New list items to add, in order, with their Word list level
# (ListLevelNumber is 1-based, matching w:ilvl + 1).
$items = @(
    @{ level = 1; text = "SQLs" },
    @{ level = 2; text = "insertTaiKhoan" },
    @{ level = 3; text = "IF NOT EXISTS ( SELECT 1 FROM taiKhoan WHERE tenDangNhap = ? OR maNhanVien = ?)" },
    @{ level = 4; text = "This sql checks if the taiKhoan/nhanVien has already existed or not. If it doesn" + [char]0x2019 + "t exists, it will perform adding the taiKhoan/nhanVien normally." }
)

$curIndex = $anchorIndex
foreach ($item in $items) {
    $curRange = $d.Paragraphs($curIndex).Range
    $curRange.InsertParagraphAfter()
    $curIndex = $curIndex + 1

    $newRange = $d.Paragraphs($curIndex).Range
    $newRange.ListFormat.ListLevelNumber = $item.level
    $newRange.Text = $item.text
}
